$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J values ---------------------------------------------------
# Rows whose own row has a row-level default style (customFormat) will pick
# that style up automatically when a bare value is written to the new,
# previously-empty J cell (mirrors real Excel's "format follows row" rule).
$ws.Range("J4").Value  = 2021
$ws.Range("J5").Value  = 5356.3
$ws.Range("J6").Value  = 9.5
$ws.Range("J8").Value  = 7.9
$ws.Range("J9").Value  = 10.5
$ws.Range("J11").Value = 9.6
$ws.Range("J12").Value = 9.4
$ws.Range("J14").Value = 14.8
$ws.Range("J15").Value = 9.1
$ws.Range("J16").Value = 9.5
$ws.Range("J17").Value = 5.9

# Row 4 has its own explicit per-cell style (not a row-level default), so the
# new J4 cell needs its formatting copied from the matching I4 cell so that
# it keeps the "year" number formatting/border instead of the plain default.
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J4").Value = 2021

# --- Empty, formatted-only cells -------------------------------------------
# J3 and J27 stay empty but need the border formatting of the row extended
# into the new column, so copy formats from the neighbouring I cells.
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)

$ws.Range("I27").Copy()
$ws.Range("J27").PasteSpecial(-4122)

# --- Selection / active cell -------------------------------------------------
# The saved workbook had the cursor resting just beyond the new data range.
$ws.Range("L27").Select() | Out-Null
